{"js": "// Update the firmware/manual version numbers shown on the cover page:\n//   \u97cc\u9ad4\u7248\u672c\uff1av2.2.2 | \u6a94\u6848\u7cfb\u7d71\u7248\u672c\uff1av1.1.0\n//   -> \u97cc\u9ad4\u7248\u672c\uff1av2.3.0 | \u6a94\u6848\u7cfb\u7d71\u7248\u672c\uff1av1.2.0\n// (commit: \"Update manual of v2.3.0 firmware\")\n\nconst body = context.document.body;\n\n// Firmware version: v2.2.2 -> v2.3.0\nconst fw = body.search(\"v2.2.2\", { matchCase: true, matchWholeWord: false });\nfw.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < fw.items.length; i++) {\n  fw.items[i].insertText(\"v2.3.0\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// File-system version: v1.1.0 -> v1.2.0\nconst fs = body.search(\"v1.1.0\", { matchCase: true, matchWholeWord: false });\nfs.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < fs.items.length; i++) {\n  fs.items[i].insertText(\"v1.2.0\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the firmware/manual version numbers shown on the cover page:\n#   \u97cc\u9ad4\u7248\u672c\uff1av2.2.2 | \u6a94\u6848\u7cfb\u7d71\u7248\u672c\uff1av1.1.0\n#   -> \u97cc\u9ad4\u7248\u672c\uff1av2.3.0 | \u6a94\u6848\u7cfb\u7d71\u7248\u672c\uff1av1.2.0\n# (commit: \"Update manual of v2.3.0 firmware\")\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2, wdFindContinue = 1\n\n# Firmware version: v2.2.2 -> v2.3.0\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"v2.2.2\", $true, $false, $false, $false, $false, $true, 1, $false, \"v2.3.0\", 2)\n\n# File-system version: v1.1.0 -> v1.2.0\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"v1.1.0\", $true, $false, $false, $false, $false, $true, 1, $false, \"v1.2.0\", 2)\n"}
